# Banana King review: drop the old "Meta description" paragraph near the
# top, and turn the closing "Prompt: ..." paragraph into two paragraphs -
# a bold title line followed by the (formerly "Meta description") blurb,
# now in italics.

$d = $word.ActiveDocument

# 1) Remove the whole "Meta description: ..." paragraph (2nd paragraph of
#    the document), including its paragraph mark, so the following
#    paragraph slides up into its place.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete() | Out-Null

# 2) The final paragraph in the document currently holds the long
#    "Prompt: Create a feature image for ..." text in italics. Insert a
#    new, empty paragraph right before it to hold the bold title line.
$lastIndex = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($lastIndex)
$lastPara.Range.InsertParagraphBefore() | Out-Null

# 3) Populate the freshly inserted (still empty) paragraph with the bold
#    title text, via a literal OOXML fragment so the run layout matches
#    the rest of the document exactly (a leading empty run followed by
#    the formatted text run).
$newTitlePara = $d.Paragraphs.Item($lastIndex)
$newTitleRange = $d.Range($newTitlePara.Range.Start, $newTitlePara.Range.End)
$titleXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Banana King slot for free: Review and Gameplay Features</w:t></w:r></w:p>'
$newTitleRange.InsertXML($titleXml) | Out-Null

# 4) Replace the text of the (still-italic) final paragraph with the old
#    meta-description blurb, leaving its existing italic formatting as-is.
$finalPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$finalRange = $d.Range($finalPara.Range.Start, $finalPara.Range.End)
$finalRange.Text = "Explore the lush forest and win big with Banana King slot. Check out our review and play Banana King for free."
